$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin / Link (B, C) are plain text already; Price / Volume (D, E) are
# numeric-looking strings ("245.08", "-0.58%") that must stay literal
# text, so force those specific cells to Text format before assigning,
# otherwise Excel auto-converts them into numbers/percentages.
# (Note: the engine only honors the first area of a comma-separated
# union Range, so NumberFormat is applied per contiguous D:E pair.)

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = '245.08'
$ws.Range("E2").Value = '-0.58%'

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = '26.98'
$ws.Range("E3").Value = '1.70%'

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = '5.070'
$ws.Range("E4").Value = '-0.31%'

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.37%'

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = '6.475'
$ws.Range("E6").Value = '-0.06%'

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8201'
$ws.Range("E7").Value = '0.78%'

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8403'
$ws.Range("E8").Value = '-0.48%'

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1329'
$ws.Range("E9").Value = '-0.99%'

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06912'
$ws.Range("E10").Value = '-0.75%'

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = '0.02860'
$ws.Range("E11").Value = '0.29%'

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09395'
$ws.Range("E12").Value = '0.03%'

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001510'
$ws.Range("E13").Value = '-1.49%'

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = '0.04091'
$ws.Range("E14").Value = '-12.07%'

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").Value = '0.006028'
$ws.Range("E15").Value = '-2.12%'

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("B16").Value = 'UpBots'
$ws.Range("C16").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D16").Value = '0.007486'
$ws.Range("E16").Value = '3,760.94%'

$ws.Range("D17:E17").NumberFormat = "@"
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '3.509'
$ws.Range("E17").Value = '-2.21%'

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").Value = '3.002'
$ws.Range("E18").Value = '-0.33%'

$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").Value = '2.313'
$ws.Range("E19").Value = '9.20%'

$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = '0.3177'
$ws.Range("E20").Value = '0.64%'

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("B21").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C21").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D21").Value = '0.03175'
$ws.Range("E21").Value = '-0.42%'

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("B22").Value = 'ProBitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D22").Value = '0.1297'
$ws.Range("E22").Value = '-1.81%'

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("B23").Value = 'MCDex'
$ws.Range("C23").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D23").Value = '3.569'
$ws.Range("E23").Value = '-5.03%'

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("B24").Value = 'ZBToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D24").Value = '0.1373'
$ws.Range("E24").Value = '1.67%'

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("B25").Value = 'One'
$ws.Range("C25").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D25").Value = '0.0005986'
$ws.Range("E25").Value = '-0.29%'

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-2.50%'

$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = '0.003964'
$ws.Range("E27").Value = '-13.63%'

$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = '0.00009793'
$ws.Range("E28").Value = '1.95%'

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03694'
$ws.Range("E40").Value = '0.69%'

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = '0.005856'
$ws.Range("E41").Value = '72.84%'

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '0.1054'
$ws.Range("E42").Value = '-22.28%'

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-10.38%'

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009379'
$ws.Range("E44").Value = '5.45%'

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005206'
$ws.Range("E45").Value = '-1.76%'

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000749'
$ws.Range("E46").Value = '-0.12%'

$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1014'
$ws.Range("E47").Value = '-15.52%'

$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002593'
$ws.Range("E48").Value = '2.92%'

$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002098'
$ws.Range("E49").Value = '-0.12%'

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001998'
$ws.Range("E50").Value = '-0.12%'
